$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header row (row 1): relabel existing columns G:K and add new columns L:O
# ---------------------------------------------------------------------------
$ws.Range("G1").Value = "점수(룰)"
$ws.Range("H1").Value = "3일상승확률(%)"
$ws.Range("I1").Value = "5일상승확률(%)"
$ws.Range("J1").Value = "10일상승확률(%)"
$ws.Range("K1").Value = "최종점수"

# Bring the new header cells (L1:O1) up to the same bold/border/centered
# style used by the rest of the header row before filling in their text.
$ws.Range("K1").Copy()
$ws.Range("L1:O1").PasteSpecial(-4122)

$ws.Range("L1").Value = "예측방식"
$ws.Range("M1").Value = "판단"
$ws.Range("N1").Value = "MACRO_SCORE"
$ws.Range("O1").Value = "MACRO_SIGNAL"

# ---------------------------------------------------------------------------
# Data rows (2:6) - refreshed for 2025-11-29, reordered and recomputed,
# plus the new 예측방식 / 판단 / MACRO_SCORE / MACRO_SIGNAL columns.
# ---------------------------------------------------------------------------
$date = "2025-11-29"
$pattern = "Pattern"
$macroScore = 85.36763896678245
$macroSignal = "🟢 완화적 (상승 우위)"
$watch = "⛔ 관망하십시오."
$buyZone = "📈 매수 관찰 구간입니다."

# Row 2: KOREA AEROSPACE
$ws.Range("B2").Value = "KOREA AEROSPACE"
$ws.Range("C2").Value = "047810.KS"
$ws.Range("D2").Value = 108900
$ws.Range("E2").ClearContents()
$ws.Range("F2").Value = 2.16
$ws.Range("G2").Value = 50
$ws.Range("H2").Value = 63
$ws.Range("I2").Value = 63
$ws.Range("J2").Value = 76
$ws.Range("K2").Value = 65.8
$ws.Range("L2").Value = $pattern
$ws.Range("M2").Value = $buyZone
$ws.Range("N2").Value = $macroScore
$ws.Range("O2").Value = $macroSignal

# Row 3: HYUNDAI ROTEM
$ws.Range("B3").Value = "HYUNDAI ROTEM"
$ws.Range("C3").Value = "064350.KS"
$ws.Range("D3").Value = 175700
$ws.Range("E3").ClearContents()
$ws.Range("F3").Value = -3.88
$ws.Range("G3").Value = 10
$ws.Range("H3").Value = 66
$ws.Range("I3").Value = 63
$ws.Range("J3").Value = 76
$ws.Range("K3").Value = 53.8
$ws.Range("L3").Value = $pattern
$ws.Range("M3").Value = $watch
$ws.Range("N3").Value = $macroScore
$ws.Range("O3").Value = $macroSignal

# Row 4: LIG Nex1
$ws.Range("B4").Value = "LIG Nex1"
$ws.Range("C4").Value = "079550.KS"
$ws.Range("D4").Value = 383000
$ws.Range("E4").ClearContents()
$ws.Range("F4").Value = -4.73
$ws.Range("G4").Value = 10
$ws.Range("H4").Value = 53
$ws.Range("I4").Value = 56
$ws.Range("J4").Value = 46
$ws.Range("K4").Value = 51
$ws.Range("L4").Value = $pattern
$ws.Range("M4").Value = $watch
$ws.Range("N4").Value = $macroScore
$ws.Range("O4").Value = $macroSignal

# Row 5: HANWHA AEROSPACE
$ws.Range("B5").Value = "HANWHA AEROSPACE"
$ws.Range("C5").Value = "012450.KS"
$ws.Range("D5").Value = 852000
$ws.Range("E5").ClearContents()
$ws.Range("F5").Value = -1.96
$ws.Range("G5").Value = 10
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 53
$ws.Range("J5").Value = 56
$ws.Range("K5").Value = 49.8
$ws.Range("L5").Value = $pattern
$ws.Range("M5").Value = $watch
$ws.Range("N5").Value = $macroScore
$ws.Range("O5").Value = $macroSignal

# Row 6: HANWHA SYSTEMS
$ws.Range("B6").Value = "HANWHA SYSTEMS"
$ws.Range("C6").Value = "272210.KS"
$ws.Range("D6").Value = 46200
$ws.Range("E6").ClearContents()
$ws.Range("F6").Value = -4.94
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 53
$ws.Range("I6").Value = 53
$ws.Range("J6").Value = 40
$ws.Range("K6").Value = 49.8
$ws.Range("L6").Value = $pattern
$ws.Range("M6").Value = $watch
$ws.Range("N6").Value = $macroScore
$ws.Range("O6").Value = $macroSignal

# ---------------------------------------------------------------------------
# Column A (날짜): every data row moves from 2025-11-28 to 2025-11-29. The
# date text must stay a literal text string (not get auto-converted into a
# date serial), so force text formatting, assign, then restore the plain
# (unstyled) look by copying the format from a never-styled neighbour cell.
# ---------------------------------------------------------------------------
$ws.Range("A2:A6").NumberFormat = "@"
$ws.Range("A2").Value = $date
$ws.Range("A3").Value = $date
$ws.Range("A4").Value = $date
$ws.Range("A5").Value = $date
$ws.Range("A6").Value = $date

$ws.Range("D2").Copy()
$ws.Range("A2:A6").PasteSpecial(-4122)

$excel.CutCopyMode = 0
